$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: establish row 81 as the template row (copy formatting from row 80) ---
$ws.Range("A80:G80").Copy($ws.Range("A81:G81"))
$ws.Range("C81").ClearContents()
$ws.Range("D81").ClearContents()
$ws.Range("A81").Font.Name = "Arial"
$ws.Range("A81").Font.Size = 10

# --- Step 2: clone row 81 formatting (incl. new Arial font) down through row 113 ---
for ($r = 82; $r -le 113; $r++) {
    $ws.Range("A80:G80").Copy($ws.Range("A" + $r + ":G" + $r))
    $ws.Range("C" + $r).ClearContents()
    $ws.Range("D" + $r).ClearContents()
    $ws.Range("A81").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Step 3: populate the UserName (A) and Email (G) values for the 33 new rows ---
$ws.Range("A81").Value = "KHPaddUser1"
$ws.Range("G81").Value = "KHPaddUser1@mailinator.com"
$ws.Range("A82").Value = "KHPaddUser2"
$ws.Range("G82").Value = "KHPaddUser2@mailinator.com"
$ws.Range("A83").Value = "KHPaddUser3"
$ws.Range("G83").Value = "KHPaddUser3@mailinator.com"
$ws.Range("A84").Value = "KHPaddUser4"
$ws.Range("G84").Value = "KHPaddUser4@mailinator.com"
$ws.Range("A85").Value = "KHPaddUser5"
$ws.Range("G85").Value = "KHPaddUser5@mailinator.com"
$ws.Range("A86").Value = "KHPaddUser6"
$ws.Range("G86").Value = "KHPaddUser6@mailinator.com"
$ws.Range("A87").Value = "SearchWhatsMarketUser1"
$ws.Range("G87").Value = "SearchWhatsMarketUser1@mailinator.com "
$ws.Range("A88").Value = "SearchWhatsMarketUser2"
$ws.Range("G88").Value = "SearchWhatsMarketUser2@mailinator.com "
$ws.Range("A89").Value = "SearchWhatsMarketUser3"
$ws.Range("G89").Value = "SearchWhatsMarketUser3@mailinator.com "
$ws.Range("A90").Value = "SearchWhatsMarketUser4"
$ws.Range("G90").Value = "SearchWhatsMarketUser4@mailinator.com "
$ws.Range("A91").Value = "SearchWhatsMarketUser5"
$ws.Range("G91").Value = "SearchWhatsMarketUser5@mailinator.com "
$ws.Range("A92").Value = "SearchWhatsMarketUser6"
$ws.Range("G92").Value = "SearchWhatsMarketUser6@mailinator.com "
$ws.Range("A93").Value = "SearchWhatsMarketUser7"
$ws.Range("G93").Value = "SearchWhatsMarketUser7@mailinator.com "
$ws.Range("A94").Value = "SearchWhatsMarketUser8"
$ws.Range("G94").Value = "SearchWhatsMarketUser8@mailinator.com "
$ws.Range("A95").Value = "SearchKnowHowUser1"
$ws.Range("G95").Value = "SearchKnowHowUser1@mailinator.com "
$ws.Range("A96").Value = "SearchKnowHowUser2"
$ws.Range("G96").Value = "SearchKnowHowUser2@mailinator.com "
$ws.Range("A97").Value = "SearchKnowHowUser3"
$ws.Range("G97").Value = "SearchKnowHowUser3@mailinator.com "
$ws.Range("A98").Value = "SearchKnowHowUser4"
$ws.Range("G98").Value = "SearchKnowHowUser4@mailinator.com "
$ws.Range("A99").Value = "SearchKnowHowUser5"
$ws.Range("G99").Value = "SearchKnowHowUser5@mailinator.com "
$ws.Range("A100").Value = "SearchKnowHowUser6"
$ws.Range("G100").Value = "SearchKnowHowUser6@mailinator.com "
$ws.Range("A101").Value = "SearchKnowHowUser7"
$ws.Range("G101").Value = "SearchKnowHowUser7@mailinator.com "
$ws.Range("A102").Value = "SearchKnowHowUser8"
$ws.Range("G102").Value = "SearchKnowHowUser8@mailinator.com "
$ws.Range("A103").Value = "AskUser1"
$ws.Range("G103").Value = "AskUser1@mailinator.com "
$ws.Range("A104").Value = "AskUser2"
$ws.Range("G104").Value = "AskUser2@mailinator.com "
$ws.Range("A105").Value = "AskUser3"
$ws.Range("G105").Value = "AskUser3@mailinator.com "
$ws.Range("A106").Value = "AskUser4"
$ws.Range("G106").Value = "AskUser4@mailinator.com "
$ws.Range("A107").Value = "AskUser5"
$ws.Range("G107").Value = "AskUser5@mailinator.com "
$ws.Range("A108").Value = "AskUser6"
$ws.Range("G108").Value = "AskUser6@mailinator.com "
$ws.Range("A109").Value = "AssetPageUser1"
$ws.Range("G109").Value = "AssetPageUser1@mailinator.com "
$ws.Range("A110").Value = "AssetPageUser2"
$ws.Range("G110").Value = "AssetPageUser2@mailinator.com "
$ws.Range("A111").Value = "AssetPageUser3"
$ws.Range("G111").Value = "AssetPageUser3@mailinator.com "
$ws.Range("A112").Value = "AssetPageUser4"
$ws.Range("G112").Value = "AssetPageUser4@mailinator.com "
$ws.Range("A113").Value = "AssetPageUser5"
$ws.Range("G113").Value = "AssetPageUser5@mailinator.com "

# --- Step 4: add mailto hyperlinks for the Email column (skip rows 81 and 112) ---
$ws.Hyperlinks.Add($ws.Range("G82"), "mailto:KHPaddUser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G83"), "mailto:KHPaddUser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G84"), "mailto:KHPaddUser4@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G85"), "mailto:KHPaddUser5@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G86"), "mailto:KHPaddUser6@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G87"), "mailto:SearchWhatsMarketUser1@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G88"), "mailto:SearchWhatsMarketUser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G89"), "mailto:SearchWhatsMarketUser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G90"), "mailto:SearchWhatsMarketUser4@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G91"), "mailto:SearchWhatsMarketUser5@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G92"), "mailto:SearchWhatsMarketUser6@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G93"), "mailto:SearchWhatsMarketUser7@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G94"), "mailto:SearchWhatsMarketUser8@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G95"), "mailto:SearchKnowHowUser1@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G96"), "mailto:SearchKnowHowUser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G97"), "mailto:SearchKnowHowUser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G98"), "mailto:SearchKnowHowUser4@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G99"), "mailto:SearchKnowHowUser5@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G100"), "mailto:SearchKnowHowUser6@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G101"), "mailto:SearchKnowHowUser7@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G102"), "mailto:SearchKnowHowUser8@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G103"), "mailto:AskUser1@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G104"), "mailto:AskUser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G105"), "mailto:AskUser3@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G106"), "mailto:AskUser4@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G107"), "mailto:AskUser5@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G108"), "mailto:AskUser6@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G109"), "mailto:AssetPageUser1@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G110"), "mailto:AssetPageUser2@mailinator.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G111"), "mailto:AssetPageUser3@mailinator.com", "", "", "AssetPageUser1@mailinator.com ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G113"), "mailto:AssetPageUser5@mailinator.com") | Out-Null

# --- Step 5: update the sheet selection to match the newly added block ---
$ws.Range("A81:G113").Select()

Write-Host "Added 33 new users (rows 81-113)"
